# Add a "time_taken" column (F) with timestamps recording how long each
# PanelApp row took to process, for metadata purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1: label + same formatting as the other header cells ---
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# --- Data cells F2:F53: one inline timestamp string per row ---
$values = @(
    "2021-10-05 10:52:47.419746",
    "2021-10-05 10:52:47.419757",
    "2021-10-05 10:52:47.419761",
    "2021-10-05 10:52:47.419763",
    "2021-10-05 10:52:47.419766",
    "2021-10-05 10:52:47.419769",
    "2021-10-05 10:52:47.419771",
    "2021-10-05 10:52:47.419774",
    "2021-10-05 10:52:47.419777",
    "2021-10-05 10:52:47.419779",
    "2021-10-05 10:52:47.419782",
    "2021-10-05 10:52:47.419784",
    "2021-10-05 10:52:47.419787",
    "2021-10-05 10:52:47.419789",
    "2021-10-05 10:52:47.419792",
    "2021-10-05 10:52:47.419794",
    "2021-10-05 10:52:47.419797",
    "2021-10-05 10:52:47.419800",
    "2021-10-05 10:52:47.419802",
    "2021-10-05 10:52:47.419805",
    "2021-10-05 10:52:47.419807",
    "2021-10-05 10:52:47.419810",
    "2021-10-05 10:52:47.419812",
    "2021-10-05 10:52:47.419815",
    "2021-10-05 10:52:47.419818",
    "2021-10-05 10:52:47.419820",
    "2021-10-05 10:52:47.419823",
    "2021-10-05 10:52:47.419825",
    "2021-10-05 10:52:47.419828",
    "2021-10-05 10:52:47.419830",
    "2021-10-05 10:52:47.419833",
    "2021-10-05 10:52:47.419835",
    "2021-10-05 10:52:47.419838",
    "2021-10-05 10:52:47.419841",
    "2021-10-05 10:52:47.419844",
    "2021-10-05 10:52:47.419846",
    "2021-10-05 10:52:47.419849",
    "2021-10-05 10:52:47.419851",
    "2021-10-05 10:52:47.419854",
    "2021-10-05 10:52:47.419856",
    "2021-10-05 10:52:47.419859",
    "2021-10-05 10:52:47.419862",
    "2021-10-05 10:52:47.419864",
    "2021-10-05 10:52:47.419867",
    "2021-10-05 10:52:47.419869",
    "2021-10-05 10:52:47.419872",
    "2021-10-05 10:52:47.419874",
    "2021-10-05 10:52:47.419877",
    "2021-10-05 10:52:47.419879",
    "2021-10-05 10:52:47.419882",
    "2021-10-05 10:52:47.419884",
    "2021-10-05 10:52:47.419887"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
